# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig -- metadata refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Bump the published version number
$ws.Range("B3").Value = "6.0.0"

# Refresh the publication date
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Fill in the real publisher (was blank) and fix the mislabeled row that
# used to say "Contact / No display for ContactDetail" -- it should be the
# Jurisdiction row.
$ws.Range("B9").Value = "Alvearie Team"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# The old sheet had an accidental duplicate "Contact" row right after --
# remove it so the rest of the metadata rows shift back into place.
$ws.Rows.Item(11).Delete()
